$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update source data values (new simulation results)
$ws.Range("B8").Value = 8257985639.0299997
$ws.Range("C8").Value = 5678945.2397499997
$ws.Range("B9").Value = 4201650652.5999999
$ws.Range("C9").Value = 8650198.6746299993

# Apply Currency format to B5:B7 and C7 (match the format used elsewhere in the column)
$ws.Range("B5").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("B6").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("B7").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("C7").NumberFormat = $ws.Range("B2").NumberFormat

# Update formulas in the summary table (rows 13-15) to reference row 4 instead of row 3
# NOTE: use single-quoted (literal) strings so PowerShell does not try to
# interpolate "$F" / "$4" as variables inside the formula text.
$ws.Range("B13").Formula = '=E5-E4'
$ws.Range("C13").Formula = '=E6-E4'
$ws.Range("D13").Formula = '=SQRT(POWER($F$4,2)+POWER(F5,2))'
$ws.Range("E13").Formula = '=SQRT(POWER($F$4,2)+POWER(F6,2))'

$ws.Range("B14").Formula = '=E7-E4'
$ws.Range("C14").Formula = '=E8-E4'
$ws.Range("D14").Formula = '=SQRT(POWER($F$4,2)+POWER(F7,2))'
$ws.Range("E14").Formula = '=SQRT(POWER($F$4,2)+POWER(F8,2))'

$ws.Range("B15").Formula = '=E9-E4'
$ws.Range("C15").Formula = '=E10-E4'
$ws.Range("D15").Formula = '=SQRT(POWER($F$4,2)+POWER(F9,2))'
$ws.Range("E15").Formula = '=SQRT(POWER($F$4,2)+POWER(F10,2))'

# Update the selected cell to B15
$ws.Range("B15").Select() | Out-Null
